$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1118.5
$ws.Range("I2").Value = 1249.6666
$ws.Range("J2").Value = 725
$ws.Range("K2").Value = 1249.6666
$ws.Range("L2").Value = 725
$ws.Range("M2").Value = -1136.6666
$ws.Range("N2").Value = -951
# Row 28
$ws.Range("H28").Value = 2024
$ws.Range("J28").Value = 3338.5
$ws.Range("L28").Value = 3338.5
$ws.Range("N28").Value = -4308.5
# Row 40
$ws.Range("H40").Value = 2900
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 3800
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 3800
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -4150
# Row 116
$ws.Range("H116").Value = 958
$ws.Range("I116").Value = 958
$ws.Range("K116").Value = 958
$ws.Range("M116").Value = 2484
# Row 134
$ws.Range("H134").Value = 75000
$ws.Range("J134").Value = 75000
$ws.Range("L134").Value = 75000
$ws.Range("N134").Value = -85140
# Row 135
$ws.Range("H135").Value = 1516
$ws.Range("I135").Value = 1516
$ws.Range("K135").Value = 13644
$ws.Range("M135").Value = -11109

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 24155.777
$ws.Range("I107").Value = 28407.8
$ws.Range("K107").Value = 28407.8
$ws.Range("M107").Value = -26487.8
# Row 134
$ws.Range("H134").Value = 3273.75
$ws.Range("I134").Value = 3273.75
$ws.Range("K134").Value = 9821.25
$ws.Range("M134").Value = -7286.25

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 1038.125
$ws.Range("I19").Value = 51.666668
$ws.Range("J19").Value = 1630
$ws.Range("K19").Value = 51.666668
$ws.Range("L19").Value = 1630
$ws.Range("M19").Value = 118.333332
$ws.Range("N19").Value = -1970
# Row 24
$ws.Range("H24").Value = 1038.125
$ws.Range("I24").Value = 51.666668
$ws.Range("J24").Value = 1630
$ws.Range("K24").Value = 51.666668
$ws.Range("L24").Value = 1630
$ws.Range("M24").Value = 118.333332
$ws.Range("N24").Value = -1970
# Row 25
$ws.Range("H25").Value = 7530.8
$ws.Range("I25").Value = 327
$ws.Range("K25").Value = 327
$ws.Range("M25").Value = -153
# Row 58
$ws.Range("H58").Value = 994.5
$ws.Range("I58").Value = 994.5
$ws.Range("K58").Value = 994.5
$ws.Range("M58").Value = -791.5
# Row 136
$ws.Range("H136").Value = 994.5
$ws.Range("I136").Value = 994.5
$ws.Range("K136").Value = 2983.5
$ws.Range("M136").Value = -433.5

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 400
$ws.Range("J23").Value = 300
$ws.Range("L23").Value = 900
$ws.Range("N23").Value = -1370
# Row 32
$ws.Range("H32").Value = 8747.5
$ws.Range("J32").Value = 8747.5
$ws.Range("L32").Value = 26242.5
$ws.Range("N32").Value = -26808.5
# Row 34
$ws.Range("H34").Value = 2877.6365
$ws.Range("I34").Value = 913.75
$ws.Range("K34").Value = 2741.25
$ws.Range("M34").Value = -2657.25
# Row 50
$ws.Range("H50").Value = 354.33334
$ws.Range("I50").Value = 354.33334
$ws.Range("K50").Value = 1063.00002
$ws.Range("M50").Value = -582.0000199999999
# Row 53
$ws.Range("H53").Value = 354.33334
$ws.Range("I53").Value = 354.33334
$ws.Range("K53").Value = 1063.00002
$ws.Range("M53").Value = -582.0000199999999
# Row 98
$ws.Range("H98").Value = 1132.5714
$ws.Range("J98").Value = 1732.75
$ws.Range("L98").Value = 5198.25
$ws.Range("N98").Value = -8194.25
# Row 140
$ws.Range("H140").Value = 435.875
$ws.Range("I140").Value = 435.875
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 1307.625
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 3872.375
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 380000000
$ws.Range("I35").Value = 380000000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 380000000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -379999702
$ws.Range("N35").ClearContents()
# Row 49
$ws.Range("H49").Value = 32092.475
$ws.Range("I49").Value = 32500
$ws.Range("J49").Value = 29919
$ws.Range("K49").Value = 32500
$ws.Range("L49").Value = 29919
$ws.Range("M49").Value = -32316
$ws.Range("N49").Value = -30287
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 2629.2
$ws.Range("I10").Value = 1286.5
$ws.Range("J10").Value = 8000
$ws.Range("K10").Value = 1286.5
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = -1146.5
$ws.Range("N10").Value = -8280
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
# Row 50
$ws.Range("H50").Value = 42666.668
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 42666.668
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 42666.668
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -43940.668
# Row 122
$ws.Range("H122").Value = 5450.7
$ws.Range("I122").Value = 4700.25
$ws.Range("K122").Value = 14100.75
$ws.Range("M122").Value = -11650.75
# Row 132
$ws.Range("H132").Value = 3398.9
$ws.Range("I132").Value = 3398.9
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10196.7
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7666.700000000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2337.25
$ws.Range("J81").Value = 1667.3334
$ws.Range("L81").Value = 3334.6668
$ws.Range("N81").Value = -5456.6668
# Row 84
$ws.Range("H84").Value = 2337.25
$ws.Range("J84").Value = 1667.3334
$ws.Range("L84").Value = 16673.334
$ws.Range("N84").Value = -27281.334
# Row 132
$ws.Range("H132").Value = 764.25
$ws.Range("I132").Value = 764.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2292.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 237.25
$ws.Range("N132").ClearContents()
